# Applies the OOXML diff to the document:
#  1. Splits the "Resurrection Appearances of our Lord before His Ascension"
#     paragraph into an empty paragraph followed by a new paragraph whose
#     run now starts with a <w:lastRenderedPageBreak/>.
#  2. Removes the <w:lastRenderedPageBreak/> that used to sit in front of
#     "Our Lord appeared to Mary Magdalene..." (it moved to the heading
#     above, per change 1).
#  3. Merges the two runs of the "...(1 Corinthians " / "15:6-7, NASB)"
#     verse into a single run and drops the <w:lastRenderedPageBreak/>
#     that used to separate them.
#  4. Splits the "BEHOLD, HE IS COMING WITH THE CLOUDS..." run into two
#     runs with a new <w:lastRenderedPageBreak/> in between.

$d = $word.ActiveDocument

function New-OpenXmlPackage([string]$bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-ParagraphByText([string]$searchText, [string]$bodyFragment) {
    $range = $d.Content
    $range.Collapse(1)
    $found = $range.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find paragraph containing: $searchText"
    }
    $p = $range.Paragraphs(1)
    $pRange = $p.Range
    $pRange.InsertXML((New-OpenXmlPackage $bodyFragment))
}

# --- Change 1 -------------------------------------------------------------
$change1 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/>' +
    '<w:t>Resurrection Appearances of our Lord before His Ascension</w:t></w:r></w:p>'
Replace-ParagraphByText "Resurrection Appearances of our Lord before His Ascension" $change1

# --- Change 2 -------------------------------------------------------------
$change2 = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t>Our Lord appeared to Mary Magdalene. She had more doctrine in her soul than any other woman in her time. Mark 16:9-11.</w:t></w:r></w:p>'
Replace-ParagraphByText "Our Lord appeared to Mary Magdalene. She had more doctrine in her soul than any other woman in her time. Mark 16:9-11." $change2

# --- Change 3 -------------------------------------------------------------
$change3 = '<w:p><w:pPr><w:pStyle w:val="Verses"/></w:pPr>' +
    '<w:r><w:t>“After that He appeared to more than five hundred brethren at one time, most of whom remain until now, but some have fallen asleep; then He appeared to James, then to all the apostles;” (1 Corinthians 15:6-7, NASB)</w:t></w:r></w:p>'
Replace-ParagraphByText "After that He appeared to more than five hundred brethren at one time" $change3

# --- Change 4 -------------------------------------------------------------
$change4 = '<w:p><w:pPr><w:pStyle w:val="Verses"/></w:pPr>' +
    '<w:r><w:t>“</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">BEHOLD, HE IS COMING WITH THE CLOUDS, and every eye will see Him, even those who pierced </w:t></w:r>' +
    '<w:r><w:lastRenderedPageBreak/><w:t>Him; and all the tribes of the earth will mourn over Him. So it is to be. Amen.</w:t></w:r>' +
    '<w:r><w:t>”</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  (Revelation 1:7, NASB)</w:t></w:r></w:p>'
Replace-ParagraphByText "BEHOLD, HE IS COMING WITH THE CLOUDS" $change4

Write-Output "All changes applied"
